$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray _GoBack bookmark that sits in the empty
#    paragraph right after "--from / --to  filters."
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Insert a new bulleted paragraph right after the paragraph that
#    ends with "when they are used together." (and before the line
#    break run that used to trail it), containing the new note about
#    Databot / AsynchronousCsvLineWriter / CSVFormatter, with a fresh
#    _GoBack bookmark placed between "CSVFormatter" and ". Reconcile."
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("when they are used together.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text 'when they are used together.'"
}

# Remember where the new paragraph break needs to go (right after the
# anchor text, before its trailing line-break run).
$splitPos = $rng.End
$rng.Collapse(0)

# --- Run 1: "Databot has " (bold) ---
$rng.InsertAfter("Databot has ")
$rng.Font.Bold = $true
$rng.Font.Name = "Garamond"
$rng.Font.Size = 10
$rng.Collapse(0)

# --- Run 2: "AsynchronousCsvLineWriter" (bold) ---
$rng.InsertAfter("AsynchronousCsvLineWriter")
$rng.Font.Bold = $true
$rng.Font.Name = "Garamond"
$rng.Font.Size = 10
$rng.Collapse(0)

# --- Run 3: ", which also uses " (not bold) ---
$rng.InsertAfter(", which also uses ")
$rng.Font.Bold = $false
$rng.Font.Name = "Garamond"
$rng.Font.Size = 10
$rng.Collapse(0)

# --- Run 4: "CSVFormatter" (not bold) ---
$rng.InsertAfter("CSVFormatter")
$rng.Font.Bold = $false
$rng.Font.Name = "Garamond"
$rng.Font.Size = 10
$rng.Collapse(0)

# Remember the position between "CSVFormatter" and ". Reconcile." for
# the new _GoBack bookmark, added once all text exists.
$bookmarkPos = $rng.Start

# --- Run 5: ". Reconcile." (not bold) ---
$rng.InsertAfter(". Reconcile.")
$rng.Font.Bold = $false
$rng.Font.Name = "Garamond"
$rng.Font.Size = 10

# Add the _GoBack bookmark at the saved (now-stable) position, using a
# freshly constructed range so earlier insertions don't drag it along.
$bmRng = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

# Finally, split the paragraph at the remembered position so that all
# the new text becomes its own list paragraph (matching the bold /
# sz20 / numId=1 formatting already used by the paragraph above), and
# the original trailing line-break run ends up inside it too.
$breakRng = $d.Range($splitPos, $splitPos)
$breakRng.InsertParagraphAfter()
